# "Adding Binary search to Insertion Sort" — mark the Kadane's Algorithm,
# Missing Number in Array and Trapping rain water rows (I3:I5) as checked on
# site ("Yes"), widen the now-meaningful "Check on site" column, and restore
# the view (zoom/selection) to what the author left it at after the edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column I ("Check on site") for rows 3-5 -> "Yes" (already "Yes" on row 2
# before this edit; rows 3-5 catch up to it).
$ws.Range("I3").Value = "Yes"
$ws.Range("I4").Value = "Yes"
$ws.Range("I5").Value = "Yes"

# Give column I (the "Check on site" column) its own explicit width now that
# every visible row in it carries a value, instead of sharing the sheet's
# default width with every column after it.
$ws.Columns.Item(9).ColumnWidth = 14.3

# View state: zoom to 125% and leave the selection on F7 (matches the
# sheetView/selection captured in the saved workbook).
$excel.ActiveWindow.Zoom = 125
$ws.Range("F7").Select()
